$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new dSF (column F) value, per the re-pulled data / mean calc fix.
$updates = @{
    2  = -1
    3  = -1
    5  = 3
    9  = 1
    13 = 0
    15 = 1
    20 = 2
    21 = -2
    25 = -2
    29 = -3
    34 = -3
    35 = -1
    42 = -2
    44 = -5
    46 = -6
    47 = -12
    48 = -5
    51 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
